$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-24 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-25 Thursday", 2) | Out-Null
$d.Content.Find.Execute("792÷2=396, 0", $true, $false, $false, $false, $false, $true, 1, $false, "386÷2=193, 0", 2) | Out-Null
$d.Content.Find.Execute("921÷5=184, 1", $true, $false, $false, $false, $false, $true, 1, $false, "735÷3=245, 0", 2) | Out-Null
$d.Content.Find.Execute("976÷7=139, 3", $true, $false, $false, $false, $false, $true, 1, $false, "681÷9=75, 6", 2) | Out-Null
$d.Content.Find.Execute("584÷5=116, 4", $true, $false, $false, $false, $false, $true, 1, $false, "804÷6=134, 0", 2) | Out-Null
$d.Content.Find.Execute("198÷5=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "998÷9=110, 8", 2) | Out-Null
$d.Content.Find.Execute("968÷3=322, 2", $true, $false, $false, $false, $false, $true, 1, $false, "190÷4=47, 2", 2) | Out-Null
$d.Content.Find.Execute("358÷9=39, 7", $true, $false, $false, $false, $false, $true, 1, $false, "759÷6=126, 3", 2) | Out-Null
$d.Content.Find.Execute("500÷8=62, 4", $true, $false, $false, $false, $false, $true, 1, $false, "116÷8=14, 4", 2) | Out-Null
$d.Content.Find.Execute("634÷6=105, 4", $true, $false, $false, $false, $false, $true, 1, $false, "148÷4=37, 0", 2) | Out-Null
$d.Content.Find.Execute("457÷7=65, 2", $true, $false, $false, $false, $false, $true, 1, $false, "869÷9=96, 5", 2) | Out-Null
$d.Content.Find.Execute("410÷8=51, 2", $true, $false, $false, $false, $false, $true, 1, $false, "289÷4=72, 1", 2) | Out-Null
$d.Content.Find.Execute("355÷4=88, 3", $true, $false, $false, $false, $false, $true, 1, $false, "847÷5=169, 2", 2) | Out-Null
$d.Content.Find.Execute("769÷6=128, 1", $true, $false, $false, $false, $false, $true, 1, $false, "170÷9=18, 8", 2) | Out-Null
$d.Content.Find.Execute("422÷9=46, 8", $true, $false, $false, $false, $false, $true, 1, $false, "102÷2=51, 0", 2) | Out-Null
$d.Content.Find.Execute("652÷6=108, 4", $true, $false, $false, $false, $false, $true, 1, $false, "618÷3=206, 0", 2) | Out-Null
$d.Content.Find.Execute("393÷4=98, 1", $true, $false, $false, $false, $false, $true, 1, $false, "338÷9=37, 5", 2) | Out-Null
$d.Content.Find.Execute("478÷3=159, 1", $true, $false, $false, $false, $false, $true, 1, $false, "628÷9=69, 7", 2) | Out-Null
$d.Content.Find.Execute("561÷9=62, 3", $true, $false, $false, $false, $false, $true, 1, $false, "568÷2=284, 0", 2) | Out-Null
$d.Content.Find.Execute("472÷6=78, 4", $true, $false, $false, $false, $false, $true, 1, $false, "930÷5=186, 0", 2) | Out-Null
$d.Content.Find.Execute("663÷9=73, 6", $true, $false, $false, $false, $false, $true, 1, $false, "954÷8=119, 2", 2) | Out-Null
$d.Content.Find.Execute("897÷7=128, 1", $true, $false, $false, $false, $false, $true, 1, $false, "199÷8=24, 7", 2) | Out-Null
$d.Content.Find.Execute("772÷6=128, 4", $true, $false, $false, $false, $false, $true, 1, $false, "767÷9=85, 2", 2) | Out-Null
$d.Content.Find.Execute("874÷6=145, 4", $true, $false, $false, $false, $false, $true, 1, $false, "555÷7=79, 2", 2) | Out-Null
$d.Content.Find.Execute("825÷3=275, 0", $true, $false, $false, $false, $false, $true, 1, $false, "759÷7=108, 3", 2) | Out-Null
$d.Content.Find.Execute("691÷3=230, 1", $true, $false, $false, $false, $false, $true, 1, $false, "124÷3=41, 1", 2) | Out-Null
